# Weekly update: insert a new price-report row for "Ají" at Vega Monumental
# Concepción. This pushes the existing rows 75-86 down to 76-87 (hence the
# new sheet dimension A1:R87) and populates the newly inserted row 75 with
# the latest week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 75, shifting rows 75:86 -> 76:87
$ws.Rows.Item(75).Insert()

# Populate the new row 75 with this week's record
$ws.Range("A75").Value = 11
$ws.Range("B75").Value = 'Vega Monumental Concepción'
$ws.Range("C75").Value = 'Bíobío'
$ws.Range("D75").Value = 44637
$ws.Range("E75").Value = 8
$ws.Range("F75").Value = 100112021
$ws.Range("G75").Value = 'Ají'
$ws.Range("H75").Value = 'Chilena(o)'
$ws.Range("I75").Value = 'Primera'
$ws.Range("J75").Value = 140
$ws.Range("K75").Value = 25000
$ws.Range("L75").Value = 26000
$ws.Range("M75").Value = 25571
$ws.Range("N75").Value = '$/saco 25 kilos'
$ws.Range("O75").Value = 'Región Metropolitana'
$ws.Range("P75").Value = 1023
$ws.Range("Q75").Value = 25
$ws.Range("R75").Value = 'Hortaliza'

# Match the date number-format used by the other rows in column D
$ws.Range("D75").NumberFormat = $ws.Range("D76").NumberFormat
